# Bladen Community College Organizations - rework column layout:
#  - swap "Organization Name" / "Categories" (now "Category") columns A/B
#  - rename several link-style headers
#  - add a new "Tiktok Link" column M
#  - widen columns G:M and re-balance A:B widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A / B data swap (category now leads, org name follows) ---
$orgNames = @(
    "Student Government Association",
    "Phi Theta Kappa Honor Society",
    "Student Volunteer Club",
    "Future Teachers Association",
    "Business Club",
    "Art Club",
    "Intramural Sports",
    "International Student Association",
    "Community Service Club"
)
$categories = @(
    "Student Government",
    "Academic",
    "Service",
    "Academic",
    "Professional",
    "Arts",
    "Athletics",
    "Cultural",
    "Service"
)

for ($i = 0; $i -lt $orgNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $categories[$i]
    $ws.Cells.Item($row, 2).Value = $orgNames[$i]
}

# --- Header row relabeling ---
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"

# --- New "Tiktok Link" column, carrying the same header styling as L1 ---
$ws.Range("L1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "Tiktok Link"

# --- Column widths (ColumnWidth is character-based; Excel stores the
#     OOXML width with a fixed 5/6-character padding offset added on
#     top, so subtract it here to land on the exact target widths) ---
$padding = 5 / 6
$widths = @{
    1  = 20
    2  = 35
    3  = 27
    4  = 11
    5  = 50
    6  = 7
    7  = 14
    8  = 15
    9  = 16
    10 = 15
    11 = 14
    12 = 14
    13 = 13
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $padding
}
